$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended to the price history table (2025-10-20).
# The leading apostrophe forces the date-like text to stay a literal
# string (matching the existing rows, which store dates as plain text,
# not Excel date serials); ClearFormats() then drops the "quote prefix"
# cell style Excel applies for that trick, so no stray formatting is
# left behind on the new cell.
$ws.Range("A66").Value = "'2025-10-20"
$ws.Range("A66").ClearFormats()

$ws.Range("B66").Value = 53.09000015258789
$ws.Range("C66").Value = 399.75
$ws.Range("D66").Value = 338.1499938964844
